# Maestro.xlsx — "Add files via upload"
# Inserts a new "Cerveza" (beer) article as the new row 4 in the
# "Artículos" sheet, pushing every existing article row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artículos")

# --- 1. Make room for the new record right after the first two articles ---
$ws.Rows.Item(4).Insert()

# --- 2. Populate the new row with the new article's data ---
$ws.Range("A4").Value = 7793147572822
$ws.Range("B4").Value = "Cerveza"
$ws.Range("C4").Value = "sin alcohol"
$ws.Range("D4").Value = "golden"
$ws.Range("E4").Value = "Imperial"
$ws.Range("F4").Value = 355
$ws.Range("G4").Value = "cm3."
$ws.Range("H4").Value = "lata"
$ws.Range("I4").Value = "Cervezas"
$ws.Range("J4").Value = "Argentina"
$ws.Range("K4").Value = 24
$ws.Range("L4").Value = $False
$ws.Range("M4").Value = $True
$ws.Range("N4").Value = "C:\EditaSoft\Imágenes de artículos\7793147572822.png"
$ws.Range("O4").Value = $True
$ws.Range("P4").Value = $True

# --- 3. The authoring tool re-saved the sheet with a fresh internal
#        sheetId (14 -> 16). Replicate that by cloning the edited sheet
#        twice (each clone consumes the next free sheetId), dropping the
#        original and the intermediate clone, and keeping only the final
#        clone — renamed back to "Artículos" in the original position. ---
$null = $ws.Copy($null, $ws)
$null = $wb.Worksheets.Item("Artículos (2)").Copy($null, $ws)

$null = $ws.Delete()
$null = $wb.Worksheets.Item("Artículos (2)").Delete()

$final = $wb.Worksheets.Item("Artículos (2) (2)")
$final.Name = "Artículos"
$null = $final.Move($wb.Worksheets.Item(1))
$null = $final.Activate()
$null = $final.Select()
